$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "75×97=" "66×70="
Replace-Text "93×79=" "82×52="
Replace-Text "61×80=" "32×49="
Replace-Text "38×38=" "26×60="
Replace-Text "63×41=" "54×25="
Replace-Text "94×92=" "44×58="
Replace-Text "84×69=" "16×87="
Replace-Text "99×58=" "81×38="
Replace-Text "68×46=" "91×86="
Replace-Text "62×19=" "96×19="
Replace-Text "85×34=" "38×87="
Replace-Text "27×75=" "15×92="
Replace-Text "43×90=" "77×64="
Replace-Text "76×98=" "13×63="
Replace-Text "31×75=" "41×59="
Replace-Text "96×73=" "63×40="
Replace-Text "68×30=" "20×40="
Replace-Text "34×41=" "80×78="
Replace-Text "86×15=" "22×89="
Replace-Text "17×56=" "17×13="
Replace-Text "12×17=" "53×38="
Replace-Text "73×65=" "69×37="
Replace-Text "31×95=" "78×45="
Replace-Text "99×70=" "12×40="
Replace-Text "41×75=" "40×39="

Write-Host "Done"
